$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Diebold-Mariano statistic (column C) and P-Value (column D) results
# after correction described in commit "Correcion a Diebold Mariano y revision de Cap1"

$ws.Range("C2").Value = -0.43819395801258
$ws.Range("D2").Value = 0.6640172576515282

$ws.Range("C3").Value = -0.8881633143687443
$ws.Range("D3").Value = 0.380693328648146

$ws.Range("C4").Value = -0.754760689904277
$ws.Range("D4").Value = 0.4555890788561443

$ws.Range("C5").Value = -0.5518122247683077
$ws.Range("D5").Value = 0.5846860333295036

$ws.Range("C6").Value = -0.8828866639496744
$ws.Range("D6").Value = 0.3834975420828419

$ws.Range("C7").Value = -0.6404550776261693
$ws.Range("D7").Value = 0.526169938700854

$ws.Range("C8").Value = -0.5110706562835502
$ws.Range("D8").Value = 0.6126045082287481

$ws.Range("C9").Value = 0.2914649982419838
$ws.Range("D9").Value = 0.7724665898568619

$ws.Range("C10").Value = 0.6916551167330753
$ws.Range("D10").Value = 0.4938502324533329

$ws.Range("C11").Value = -0.003502609520207065
$ws.Range("D11").Value = 0.9972257985826867

$wb.Save()
